$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.274.46'
Set-TextValue $ws.Range("E2") '  +0.12%  '

Set-TextValue $ws.Range("D3") '1.679.34'

Set-TextValue $ws.Range("D4") '1.007'
Set-TextValue $ws.Range("E4") '  +0.15%  '

Set-TextValue $ws.Range("D5") '218.41'
Set-TextValue $ws.Range("E5") '  +0.36%  '

Set-TextValue $ws.Range("D6") '0.5256'
Set-TextValue $ws.Range("E6") '  +2.36%  '

Set-TextValue $ws.Range("E7") '  +0.17%  '

Set-TextValue $ws.Range("D8") '0.2698'
Set-TextValue $ws.Range("E8") '  +1.45%  '

Set-TextValue $ws.Range("D9") '0.06431'
Set-TextValue $ws.Range("E9") '  +0.73%  '

Set-TextValue $ws.Range("D10") '21.99'
Set-TextValue $ws.Range("E10") '  +1.89%  '

Set-TextValue $ws.Range("D11") '0.07513'
Set-TextValue $ws.Range("E11") '  +1.68%  '

Set-TextValue $ws.Range("D12") '1.703.08'
Set-TextValue $ws.Range("E12") '  +1.53%  '

Set-TextValue $ws.Range("D13") '4.541'
Set-TextValue $ws.Range("E13") '  -0.26%  '

Set-TextValue $ws.Range("D14") '0.5798'
Set-TextValue $ws.Range("E14") '  -0.76%  '

Set-TextValue $ws.Range("D15") '0.000008478'
Set-TextValue $ws.Range("E15") '  -2.13%  '

Set-TextValue $ws.Range("D16") '64.18'
Set-TextValue $ws.Range("E16") '  -0.73%  '

Set-TextValue $ws.Range("D17") '26.290.02'
Set-TextValue $ws.Range("E17") '  -0.09%  '

Set-TextValue $ws.Range("D18") '4.920'
Set-TextValue $ws.Range("E18") '  -0.93%  '

Set-TextValue $ws.Range("D19") '1.007'
Set-TextValue $ws.Range("E19") '  +0.12%  '

Set-TextValue $ws.Range("D20") '10.86'
Set-TextValue $ws.Range("E20") '  -0.26%  '

Set-TextValue $ws.Range("D21") '189.28'

Set-TextValue $ws.Range("D22") '6.193'
Set-TextValue $ws.Range("E22") '  -0.42%  '

Set-TextValue $ws.Range("D23") '1.008'
Set-TextValue $ws.Range("E23") '  +0.12%  '

Set-TextValue $ws.Range("D24") '144.91'
Set-TextValue $ws.Range("E24") '  +0.64%  '

Set-TextValue $ws.Range("D25") '7.717'
Set-TextValue $ws.Range("E25") '  +0.82%  '

Set-TextValue $ws.Range("D26") '0.1238'
Set-TextValue $ws.Range("E26") '  +4.67%  '

Set-TextValue $ws.Range("D27") '15.81'
Set-TextValue $ws.Range("E27") '  +0.96%  '

Set-TextValue $ws.Range("D28") '0.06575'
Set-TextValue $ws.Range("E28") '  +9.77%  '

Set-TextValue $ws.Range("E29") '  +5.74%  '

Set-TextValue $ws.Range("D30") '1.328'
Set-TextValue $ws.Range("E30") '  +0.18%  '

Set-TextValue $ws.Range("D31") '3.575'
Set-TextValue $ws.Range("E31") '  +1.32%  '

Set-TextValue $ws.Range("D32") '3.566'
Set-TextValue $ws.Range("E32") '  +0.99%  '

Set-TextValue $ws.Range("E33") '  +0.88%  '

Set-TextValue $ws.Range("D34") '1.024'
Set-TextValue $ws.Range("E34") '  +0.74%  '

Set-TextValue $ws.Range("D35") '0.6186'
Set-TextValue $ws.Range("E35") '  +2.62%  '

Set-TextValue $ws.Range("D36") '2.399'

Set-TextValue $ws.Range("D37") '2.704'
Set-TextValue $ws.Range("E37") '  +2.19%  '

Set-TextValue $ws.Range("D38") '6.385'
Set-TextValue $ws.Range("E38") '  +4.86%  '

Set-TextValue $ws.Range("B39") 'Maker'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D39") '1.105.06'
Set-TextValue $ws.Range("E39") '  +2.15%  '

Set-TextValue $ws.Range("B40") 'VeChain'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D40") '0.01618'
Set-TextValue $ws.Range("E40") '  +0.01%  '

Set-TextValue $ws.Range("D41") '0.8736'
Set-TextValue $ws.Range("E41") '  +0.65%  '

Set-TextValue $ws.Range("D42") '1.014'
Set-TextValue $ws.Range("E42") '  +0.35%  '

Set-TextValue $ws.Range("D43") '100.46'
Set-TextValue $ws.Range("E43") '  +0.19%  '

Set-TextValue $ws.Range("D44") '1.827.57'
Set-TextValue $ws.Range("E44") '  +0.38%  '

Set-TextValue $ws.Range("E45") '  +1.53%  '

Set-TextValue $ws.Range("D46") '56.74'
Set-TextValue $ws.Range("E46") '  +1.00%  '

Set-TextValue $ws.Range("D47") '1.007'
Set-TextValue $ws.Range("E47") '  -0.42%  '

Set-TextValue $ws.Range("D48") '8.109'
Set-TextValue $ws.Range("E48") '  +0.30%  '

Set-TextValue $ws.Range("D49") '0.05268'
Set-TextValue $ws.Range("E49") '  +0.99%  '

Set-TextValue $ws.Range("D50") '0.4300'
Set-TextValue $ws.Range("E50") '  +0.07%  '

Set-TextValue $ws.Range("D51") '6.042'
Set-TextValue $ws.Range("E51") '  +2.48%  '
